# Apply the "Updated symbol list" crypto-price refresh to Sheet1.
#
# The source cells are plain text (inline strings), including the
# numeric-looking Price column, so every value below is written with a
# leading apostrophe to force Excel to store it as literal text rather
# than auto-converting it to a Number. The apostrophe marks the cell
# with a "quote prefix" style; the final pass resets that style back to
# Normal on every touched cell so only the cell *values* change, just
# like in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.93"

$ws.Range("D3").Value = "'25.20"

$ws.Range("B4").Value = "'LEO"
$ws.Range("C4").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "'3.500"
$ws.Range("E4").Value = "'3LEOLEO"

$ws.Range("B5").Value = "'HuobiToken"
$ws.Range("C5").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "'5.173"
$ws.Range("E5").Value = "'4HuobiTokenHT"

$ws.Range("B6").Value = "'Cronos"
$ws.Range("C6").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "'0.05737"
$ws.Range("E6").Value = "'5CronosCRO"

$ws.Range("B7").Value = "'KuCoinToken"
$ws.Range("C7").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.499"
$ws.Range("E7").Value = "'6KuCoinTokenKCS"

$ws.Range("B8").Value = "'GateToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.111"
$ws.Range("E8").Value = "'7GateTokenGT"

$ws.Range("B9").Value = "'MXToken"
$ws.Range("C9").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.8094"
$ws.Range("E9").Value = "'8MXTokenMX"

$ws.Range("B10").Value = "'FTXToken"
$ws.Range("C10").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "'0.8406"
$ws.Range("E10").Value = "'9FTXTokenFTT"

$ws.Range("B11").Value = "'One"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D11").Value = "'0.009669"
$ws.Range("E11").Value = "'10OneONEBestin24h"

$ws.Range("B12").Value = "'WazirX"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1339"
$ws.Range("E12").Value = "'11WazirXWRX"

$ws.Range("B13").Value = "'MandalaExchangeToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.06952"
$ws.Range("E13").Value = "'12MandalaExchangeTokenMDX"

$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02833"
$ws.Range("E14").Value = "'13BitrueCoinBTR"

$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09367"
$ws.Range("E15").Value = "'14BitMartTokenBMX"

$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001521"
$ws.Range("E16").Value = "'15BitForexTokenBF"

$ws.Range("B17").Value = "'TigerCash"
$ws.Range("C17").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006078"
$ws.Range("E17").Value = "'16TigerCashTCH"

$ws.Range("D19").Value = "'0.3197"

$ws.Range("D20").Value = "'0.03132"

$ws.Range("D22").Value = "'3.758"

$ws.Range("D23").Value = "'0.04656"

$ws.Range("D25").Value = "'0.001237"

$ws.Range("D26").Value = "'0.004264"

$ws.Range("D27").Value = "'0.00009699"
$ws.Range("E27").Value = "'26NitroExNTX"

$ws.Range("D40").Value = "'0.03612"

$ws.Range("D41").Value = "'0.006349"

$ws.Range("D42").Value = "'0.1050"

$ws.Range("D44").Value = "'0.007342"

$ws.Range("D45").Value = "'0.00005304"

$ws.Range("D47").Value = "'0.1500"

$ws.Range("D48").Value = "'0.002286"

# Clear the quote-prefix formatting introduced by the apostrophes above
# so the edited cells keep their original (default) style.
$ws.Range("D2,D3,B4,C4,D4,E4,B5,C5,D5,E5,B6,C6,D6,E6,B7,C7,D7,E7,B8,C8,D8,E8,B9,C9,D9,E9,B10,C10,D10,E10,B11,C11,D11,E11,B12,C12,D12,E12,B13,C13,D13,E13,B14,C14,D14,E14,B15,C15,D15,E15,B16,C16,D16,E16,B17,C17,D17,E17,D19,D20,D22,D23,D25,D26,D27,E27,D40,D41,D42,D44,D45,D47,D48").Style = "Normal"
